$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new quarterly sheet "2022-Q4" right before "2022-Q3".
#    Copying the existing "2022-Q3" sheet keeps all header/row styling
#    identical (same layout as the other quarter sheets), and leaves
#    every other sheet's own worksheet part untouched - they simply
#    shift down one slot when the workbook is saved.
# ------------------------------------------------------------------
$srcQ3 = $wb.Worksheets.Item("2022-Q3")
$srcQ3.Copy($srcQ3)
$newQ4 = $wb.ActiveSheet
$newQ4.Name = "2022-Q4"

# Fill in the 2022-Q4 figures (columns D:G are stored as text, H as a number,
# matching the other quarter sheets).
$newQ4.Range("D2:G3").NumberFormat = "@"

$newQ4.Range("D2").Value = "1.18"
$newQ4.Range("E2").Value = "93.72"
$newQ4.Range("F2").Value = "1.89"
$newQ4.Range("G2").Value = "0.0223"
$newQ4.Range("H2").Value = 1

$newQ4.Range("D3").Value = "0.89"
$newQ4.Range("E3").Value = "93.72"
$newQ4.Range("F3").Value = "1.89"
$newQ4.Range("G3").Value = "0.0168"
$newQ4.Range("H3").Value = 1

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push the existing quarter rows
#    down by one and add a new row for 2022-Q4 right after the header.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A5:D5").Copy()
$total.Range("A6").PasteSpecial()
$total.Range("A4:D4").Copy()
$total.Range("A5").PasteSpecial()
$total.Range("A3:D3").Copy()
$total.Range("A4").PasteSpecial()
$total.Range("A2:D2").Copy()
$total.Range("A3").PasteSpecial()

# Restore the incrementing index in column A for the last (new) row and
# make sure it keeps the same formatting as the rest of column A.
$total.Range("A6").Value = 4
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

# New 2022-Q4 row.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04

# Leave the workbook on the same active sheet ("总计") it was on before the edit.
$total.Activate()
